$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '25.974.94'
$ws.Range("D2").Style = "Normal"

$ws.Range("E2").Value = '  +0.49%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.588.29'
$ws.Range("D3").Style = "Normal"

$ws.Range("E3").Value = '  +0.18%  '

$ws.Range("E4").Value = '  -0.19%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '210.40'
$ws.Range("D5").Style = "Normal"

$ws.Range("E5").Value = '  +0.22%  '

$ws.Range("E6").Value = '  -0.21%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.480'
$ws.Range("D7").Style = "Normal"

$ws.Range("E7").Value = '  +0.45%  '

$ws.Range("E8").Value = '  -0.40%  '

$ws.Range("E9").Value = '  -0.93%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '17.91'
$ws.Range("D10").Style = "Normal"

$ws.Range("E10").Value = '  -0.66%  '

$ws.Range("E11").Value = '  +2.14%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.810.06'
$ws.Range("D12").Style = "Normal"

$ws.Range("E12").Value = '  +0.25%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.583.44'
$ws.Range("D13").Style = "Normal"

$ws.Range("E13").Value = '  -0.10%  '

$ws.Range("E14").Value = '  -1.13%  '

$ws.Range("E15").Value = '  +0.00%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '25.970.85'
$ws.Range("D16").Style = "Normal"

$ws.Range("E16").Value = '  +0.56%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '60.06'
$ws.Range("D17").Style = "Normal"

$ws.Range("E17").Value = '  +0.59%  '

$ws.Range("E18").Value = '  -0.32%  '

$ws.Range("E19").Value = '  -0.22%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '199.80'
$ws.Range("D20").Style = "Normal"

$ws.Range("E20").Value = '  +4.32%  '

$ws.Range("E21").Value = '  +0.83%  '

$ws.Range("E22").Value = '  -2.06%  '

$ws.Range("E23").Value = '  +0.56%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.87'
$ws.Range("D24").Style = "Normal"

$ws.Range("E24").Value = '  +9.63%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '142.67'
$ws.Range("D25").Style = "Normal"

$ws.Range("E25").Value = '  +0.46%  '

$ws.Range("E26").Value = '  -0.21%  '

$ws.Range("E27").Value = '  -8.19%  '

$ws.Range("E28").Value = '  -0.31%  '

$ws.Range("E29").Value = '  -0.09%  '

$ws.Range("E30").Value = '  +0.15%  '

$ws.Range("E31").Value = '  +0.43%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.11'
$ws.Range("D32").Style = "Normal"

$ws.Range("E32").Value = '  +0.09%  '

$ws.Range("E33").Value = '  -3.10%  '

$ws.Range("E34").Value = '  -1.54%  '

$ws.Range("E35").Value = '  -0.73%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.123.85'
$ws.Range("D36").Style = "Normal"

$ws.Range("E36").Value = '  +2.19%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.0162'
$ws.Range("D37").Style = "Normal"

$ws.Range("E37").Value = '  +8.41%  '

$ws.Range("E38").Value = '  -0.21%  '

$ws.Range("E39").Value = '  -1.60%  '

$ws.Range("E40").Value = '  +0.70%  '

$ws.Range("E41").Value = '  -2.68%  '

$ws.Range("E42").Value = '  -5.35%  '

$ws.Range("B43").Value = 'FraxShare'

$ws.Range("C43").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.10'
$ws.Range("D43").Style = "Normal"

$ws.Range("E43").Value = '  -1.42%  '

$ws.Range("B44").Value = 'RocketPoolETH'

$ws.Range("C44").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.721.25'
$ws.Range("D44").Style = "Normal"

$ws.Range("E44").Value = '  +0.12%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '92.00'
$ws.Range("D45").Style = "Normal"

$ws.Range("E45").Value = '  -2.01%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.49'
$ws.Range("D46").Style = "Normal"

$ws.Range("E46").Value = '  -1.37%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '53.27'
$ws.Range("D47").Style = "Normal"

$ws.Range("E47").Value = '  +0.06%  '

$ws.Range("E48").Value = '  -1.16%  '

$ws.Range("E49").Value = '  -0.19%  '

$ws.Range("E50").Value = '  +0.08%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0₇0938'
$ws.Range("D51").Style = "Normal"

$ws.Range("E51").Value = '  -16.27%  '
